# Auto-generated edit script applying the diff to cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

# Row 2
Set-TextCell "D2" "57.755.11"
$ws.Range("E2").Value = "  -4.49%  "

# Row 3
Set-TextCell "D3" "2.955.18"
$ws.Range("E3").Value = "  -1.00%  "

# Row 4
Set-TextCell "D4" "1.00"
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
Set-TextCell "D5" "555.58"
$ws.Range("E5").Value = "  -2.64%  "

# Row 6
Set-TextCell "D6" "128.43"
$ws.Range("E6").Value = "  +3.04%  "

# Row 7
Set-TextCell "D7" "1.00"
$ws.Range("E7").Value = "  +0.32%  "

# Row 8
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextCell "D8" "2.947.36"
$ws.Range("E8").Value = "  -1.04%  "

# Row 9
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell "D9" "0.509"
$ws.Range("E9").Value = "  +2.64%  "

# Row 10
$ws.Range("E10").Value = "  -3.85%  "

# Row 11
$ws.Range("E11").Value = "  -5.31%  "

# Row 12
Set-TextCell "D12" "0.444"
$ws.Range("E12").Value = "  +2.25%  "

# Row 13
Set-TextCell "D13" "0.0000220"
$ws.Range("E13").Value = "  -1.21%  "

# Row 14
Set-TextCell "D14" "32.73"
$ws.Range("E14").Value = "  +0.80%  "

# Row 15
$ws.Range("E15").Value = "  +1.52%  "

# Row 16
Set-TextCell "D16" "3.442.20"
$ws.Range("E16").Value = "  -1.03%  "

# Row 17
Set-TextCell "D17" "2.957.83"
$ws.Range("E17").Value = "  -0.59%  "

# Row 18
Set-TextCell "D18" "6.65"
$ws.Range("E18").Value = "  +8.44%  "

# Row 19
Set-TextCell "D19" "57.838.89"
$ws.Range("E19").Value = "  -4.34%  "

# Row 20
Set-TextCell "D20" "415.02"
$ws.Range("E20").Value = "  -2.24%  "

# Row 21
Set-TextCell "D21" "13.07"
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
Set-TextCell "D22" "0.681"
$ws.Range("E22").Value = "  +3.55%  "

# Row 23
Set-TextCell "D23" "6.96"
$ws.Range("E23").Value = "  -2.23%  "

# Row 24
$ws.Range("E24").Value = "  +0.70%  "

# Row 25
Set-TextCell "D25" "79.14"
$ws.Range("E25").Value = "  +0.52%  "

# Row 26
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
Set-TextCell "D27" "1.00"
$ws.Range("E27").Value = "  +0.43%  "

# Row 28
Set-TextCell "D28" "2.50"
$ws.Range("E28").Value = "  -0.23%  "

# Row 29
Set-TextCell "D29" "7.49"
$ws.Range("E29").Value = "  +4.99%  "

# Row 30
Set-TextCell "D30" "1.97"
$ws.Range("E30").Value = "  +4.83%  "

# Row 31
Set-TextCell "D31" "6.13"
$ws.Range("E31").Value = "  +1.79%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D32" "0.103"
$ws.Range("E32").Value = "  +11.32%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D33" "24.99"
$ws.Range("E33").Value = "  -0.59%  "

# Row 34
Set-TextCell "D34" "5.59"
$ws.Range("E34").Value = "  +1.18%  "

# Row 35
Set-TextCell "D35" "0.933"
$ws.Range("E35").Value = "  -1.09%  "

# Row 36
$ws.Range("E36").Value = "  -8.00%  "

# Row 37
Set-TextCell "D37" "48.29"
$ws.Range("E37").Value = "  -2.14%  "

# Row 38
Set-TextCell "D38" "0.0₃0656"
$ws.Range("E38").Value = "  +0.98%  "

# Row 39
Set-TextCell "D39" "8.29"
$ws.Range("E39").Value = "  +5.90%  "

# Row 40
$ws.Range("E40").Value = "  +7.38%  "

# Row 41
$ws.Range("E41").Value = "  -0.44%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D42" "376.01"
$ws.Range("E42").Value = "  +0.46%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D43" "0.0345"
$ws.Range("E43").Value = "  -2.94%  "

# Row 44
Set-TextCell "D44" "2.644.50"
$ws.Range("E44").Value = "  +0.43%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D46" "0.235"
$ws.Range("E46").Value = "  +0.84%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D47" "120.53"
$ws.Range("E47").Value = "  +1.57%  "

# Row 48
Set-TextCell "D48" "0.109"
$ws.Range("E48").Value = "  +2.90%  "

# Row 49
Set-TextCell "D49" "1.98"
$ws.Range("E49").Value = "  +1.24%  "

# Row 50
Set-TextCell "D50" "23.24"
$ws.Range("E50").Value = "  +0.13%  "

# Row 51
Set-TextCell "D51" "1.99"
$ws.Range("E51").Value = "  +1.07%  "
